$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat is forced to Text ("@") before assignment so that Excel
# COM does not auto-coerce these numeric/percent-looking strings into
# actual numbers, preserving the original text representation.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.42%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "13.35%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.155"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.15%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07764"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.34%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.347"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "8.82%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.81%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.946"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.23%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9313"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.93%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09926"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "11.48%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1792"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.82%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08588"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.22%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03313"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.42%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09924"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.47%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001501"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.95%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005753"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.78%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.38%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.147"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.42%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3366"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.74%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.304"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.56%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2303"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.59%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04536"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.50%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001216"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.15%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004375"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.29%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.08%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.04%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01792"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "10.81%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04788"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.73%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007771"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.05%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1411"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.13%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.006839"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-23.48%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002074"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.59%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009447"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006115"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.17%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.06%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.995"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "33.97%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.06%"
